$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Solde structurel (en % du PIB)"
$ws.Range("B3").Value = "-4,6"
$ws.Range("C3").Value = "-3,3"

$ws.Range("A4").Value = "Solde conjoncturel (en % du PIB)"
$ws.Range("B4").Value = "-0,4"
$ws.Range("C4").Value = "-0,4"

$ws.Range("A5").Value = "Solde des mesures ponctuelles et temporaires (en % du PIB)"
$ws.Range("B5").Value = "-0,1"
$ws.Range("C5").Value = "-0,1"

$ws.Range("A6").Value = "Solde effectif (en % du PIB)"
$ws.Range("B6").Value = "-5,0"
$ws.Range("C6").Value = "-3,7"
